$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2.0
$ws.Range("A4").Value = 4.0
$ws.Range("A5").Value = 5.0
$ws.Range("A6").Value = 3.0
